$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active/selected sheet (matches workbook.xml
# activeTab + sheetView tabSelected moving from NewLoanInput to this sheet)
$ws.Activate()

# Insert a new (blank) column before column N, shifting the old N/O/P data
# one column to the right (N->O, O->P, P->Q). Excel copies the column
# width/format from the column to the left (M) onto the freshly inserted one.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Restore the recorded selection on the sheet after the edit
$ws.Range("S6").Select()
